$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the single "Search" result entry: lower-case search key (A),
# display title (B) and formatted price (C). Swap "el principito" for
# "crimen y castigo" per the updated validation data.
$ws.Range("A2").Value = "crimen y castigo"
$ws.Range("B2").Value = "CRIMEN Y CASTIGO"

# The price must stay a plain text value (e.g. "$50,000"), not a currency
# number, to match the original cell's stored type. Use single quotes so
# PowerShell doesn't treat "$50" as a variable, force the cell to Text
# before assigning so Excel doesn't auto-convert it to a number, then
# clear the formatting back off so no extra style gets stamped on the
# cell (the source cell carries no explicit style either).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = '$50,000'
$ws.Range("C2").ClearFormats()
